$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width ---
$ws.Columns("A").ColumnWidth = 40.3

# --- Row 5 ---
$ws.Range("B5").Value = "Nathaniel Branden"
$ws.Range("A5").Value = "La vera autostima è quella che proviamo per noi stessi quando qualcosa va storto."
$ws.Range("D5").Value = "Selbstwert"

# --- Row 6 (with custom fonts) ---
$ws.Range("B6").Value = "Denis Waitley"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Font.ColorIndex = 1
$ws.Range("B6").Font.Size = 10
$ws.Range("B6").Font.Italic = $true

$ws.Range("A6").Value = "Quando sei capace di applaudire te stesso, è molto più facile applaudire gli altri."
$ws.Range("C3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Font.ThemeColor = 1

$ws.Range("D6").Value = "Selbstwert, Italienisch"

# --- Row 7 ---
$ws.Range("A7").Value = "You never change things by fighting the existing reality. To change something, build a new model that makes the existing model obsolete"
$ws.Range("B7").Value = "Buckminster Fuller"
$ws.Range("D7").Value = "Transformation"

# --- Selection state ---
$ws.Range("A1:D7").Select() | Out-Null
